$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column F ("Oct 07" re-run) with the latest Jenkins build's
# loading-time numbers, mirroring the header of column E.
$ws.Range("F1").Style = "Normal"
$ws.Range("F1").Value = $ws.Range("E1").Value2

$ws.Range("F2").Style = "Normal"
$ws.Range("F2").Value = 15

$ws.Range("F3").Style = "Normal"
$ws.Range("F3").Value = 5

$ws.Range("F4").Style = "Normal"
$ws.Range("F4").Value = 2

$ws.Range("F5").Style = "Normal"
$ws.Range("F5").Value = 0

$ws.Range("F6").Style = "Normal"
$ws.Range("F6").Value = 0

$ws.Range("F7").Style = "Normal"
$ws.Range("F7").Value = 0
